$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Find the paragraph that ends the "Lighting" section (the
# "Point Light ... is." paragraph). The new "Skybox" content, plus
# the relocated "_GoBack" bookmark, is inserted right after it and
# before the document's trailing empty paragraph.
# ------------------------------------------------------------------
$lightingIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Point Light is being used*") {
        $lightingIndex = $i
    }
}

$anchorRange = $d.Paragraphs.Item($lightingIndex).Range
$anchorRange.Collapse(0)

# ------------------------------------------------------------------
# Insert the three new (still empty / unformatted) paragraphs first,
# while nothing bold exists yet to be inherited.
# ------------------------------------------------------------------
$anchorRange.InsertParagraphAfter()
$bookmarkIndex = $lightingIndex + 1

$bmPara = $d.Paragraphs.Item($bookmarkIndex)
$bmRange = $bmPara.Range
$bmRange.Collapse(0)
$bmRange.InsertParagraphAfter()
$skyboxIndex = $bookmarkIndex + 1

$sbPara = $d.Paragraphs.Item($skyboxIndex)
$sbRange = $sbPara.Range
$sbRange.Collapse(0)
$sbRange.InsertParagraphAfter()
$descIndex = $skyboxIndex + 1

# ------------------------------------------------------------------
# 3) Descriptive paragraph explaining the skybox implementation
#    (filled in before the heading above it is made bold, so it does
#    not inherit bold formatting).
# ------------------------------------------------------------------
$descPara = $d.Paragraphs.Item($descIndex)
$descPara.Range.Text = "We create an enormous 3D cube and we are putting textures inside it. We have downloaded multiple images which when placed together can resemble a real 3D environment. Each picture is carefully placed in a specific interior side of the cube. Inside this cube we will place our normal objects, and because the cube is so large it can fool the user that this is the whole sky and not just a cube side. For the picture used in the skybox consult ‘resources’ "

# ------------------------------------------------------------------
# 2) Bold "Skybox:" heading paragraph.
# ------------------------------------------------------------------
$sbPara = $d.Paragraphs.Item($skyboxIndex)
$sbPara.Range.Text = "Skybox:"
$sbPara = $d.Paragraphs.Item($skyboxIndex)
$sbPara.Range.Bold = 1

# ------------------------------------------------------------------
# 1) Paragraph that hosts the relocated "_GoBack" bookmark.
#    A bookmark cannot be anchored at a collapsed Range sitting in a
#    brand-new, totally empty paragraph, so a placeholder character
#    is typed first, the bookmark is dropped at its start, and the
#    placeholder is erased again (the bookmark tags stay behind).
#    Adding a new "_GoBack" bookmark also automatically removes the
#    document's previous one (paragraph 3), collapsing it back down
#    to a plain empty paragraph, exactly as the target edit requires.
# ------------------------------------------------------------------
$bmPara = $d.Paragraphs.Item($bookmarkIndex)
$bmPara.Range.Text = "x"
$bmPara = $d.Paragraphs.Item($bookmarkIndex)
$bmStart = $bmPara.Range.Start
$bmPoint = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bmPoint)
$bmPara = $d.Paragraphs.Item($bookmarkIndex)
$clearRange = $d.Range($bmPara.Range.Start, $bmPara.Range.End - 1)
$clearRange.Text = ""

$d.Save()
